# Commit: "Build site at 2022-09-26 16:07:08 UTC"
# The generator re-rendered 8800011.xlsx: several label/value rows shifted
# down by one, a new "Programa resumido: Semestral" row was introduced,
# the long-form Objetivos/Programa paragraphs and the Bibliografia entry
# were dropped, and the sheet shrank from A1:C24 to A1:C23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rewrite cell values for the affected rows (10-23) ---
$ws.Range("A10").Value = 'Objetivos:'
$ws.Range("B10").Value = '8452037 - Elisabeth Pinheiro da Silva Kondracki de Alcantara'
$ws.Range("C10").Value = '8452037 - Elisabeth Pinheiro da Silva Kondracki de Alcantara'

$ws.Range("A11").Value = 'Objectives:'
$ws.Range("B11").Value = '1. Approximate student of their vocal apparatus, at the level of spoken and sung words, in their individual and collective expression (choral). 2. Encourage the student to experience the choral repertoire and its role in the development of musical language. 3. Provide the student with the opportunity to vocal music together, with technical learning parameters such as tuning, precision, balance, phrasing etc.'
$ws.Range("C11").Value = '1. Approximate student of their vocal apparatus, at the level of spoken and sung words, in their individual and collective expression (choral). 2. Encourage the student to experience the choral repertoire and its role in the development of musical language. 3. Provide the student with the opportunity to vocal music together, with technical learning parameters such as tuning, precision, balance, phrasing etc.'

$ws.Range("A12").Value = 'Docentes responsáveis:'

$ws.Range("A13").Value = 'Programa resumido:'
$ws.Range("B13").Value = 'Semestral'
$ws.Range("C13").Value = 'Semestral'

$ws.Range("A14").Value = 'Short syllabus:'
$ws.Range("B14").Value = 'Resumption and improvement of repertoires and technical assumptions of previous discipline Choir Singing.Vocal classification. Breathing into the singing. Placing the emission in Bocca Chiusa. Choral singing in unison. Choral singing in canon. Choral singing in other polyphonic formations. Coral reading.'
$ws.Range("C14").Value = 'Resumption and improvement of repertoires and technical assumptions of previous discipline Choir Singing.Vocal classification. Breathing into the singing. Placing the emission in Bocca Chiusa. Choral singing in unison. Choral singing in canon. Choral singing in other polyphonic formations. Coral reading.'

$ws.Range("A15").Value = 'Programa:'
$ws.Range("B15").Value = '01/01/2017'
$ws.Range("C15").Value = '01/01/2017'

$ws.Range("A16").Value = 'Syllabus:'
$ws.Range("B16").Value = 'Resumption and improvement of repertoires and technical assumptions of previous discipline Choir Singing. Vocal classification.Breathing into the singing - Exercises for locating low and average breathing. Support and air column. Placing the emission in Bocca Chiusa. - relaxation of the mandible, tongue position, the soft palate suspension, local perception exercises for where the voice is being put, support connection and vocal emission, passing the Bocca Chiusa for vowels and other nasal and guttural sounds. Choral singing in unison. - The choral singing in unison, tuning, timbre uniformity, rhythmic precision. Choral singing in canon. Choral singing in other polyphonic formations. Coral reading.- testing of harmonic relaying and listening to 1st view the various voices, memorization, music theory basics. Assembling and improvement of pieces - promoting the application of learnt techniques. Connection between diaphragm and vocal emission.'
$ws.Range("C16").Value = 'Resumption and improvement of repertoires and technical assumptions of previous discipline Choir Singing. Vocal classification.Breathing into the singing - Exercises for locating low and average breathing. Support and air column. Placing the emission in Bocca Chiusa. - relaxation of the mandible, tongue position, the soft palate suspension, local perception exercises for where the voice is being put, support connection and vocal emission, passing the Bocca Chiusa for vowels and other nasal and guttural sounds. Choral singing in unison. - The choral singing in unison, tuning, timbre uniformity, rhythmic precision. Choral singing in canon. Choral singing in other polyphonic formations. Coral reading.- testing of harmonic relaying and listening to 1st view the various voices, memorization, music theory basics. Assembling and improvement of pieces - promoting the application of learnt techniques. Connection between diaphragm and vocal emission.'

$ws.Range("A17").Value = 'Avaliação:'
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()

$ws.Range("A18").Value = 'Método:'
$ws.Range("B18").Value = '8452037 - Elisabeth Pinheiro da Silva Kondracki de Alcantara'
$ws.Range("C18").Value = '8452037 - Elisabeth Pinheiro da Silva Kondracki de Alcantara'

$ws.Range("A19").Value = 'Critério:'
$ws.Range("B19").Value = 'A cada semestre é proposto um programa com cerca de 8 (oito) peças, sendo duas ou três de semestres anteriores e, consequentemente, cinco ou seis inéditas a ser apresentado pelo CORAL da EEL-USP em performances públicas definidas durante o período letivo.'
$ws.Range("C19").Value = 'A cada semestre é proposto um programa com cerca de 8 (oito) peças, sendo duas ou três de semestres anteriores e, consequentemente, cinco ou seis inéditas a ser apresentado pelo CORAL da EEL-USP em performances públicas definidas durante o período letivo.'

$ws.Range("A20").Value = 'Norma de recuperação:'
$ws.Range("B20").Value = 'Sendo uma atividade prática e de grupo, fica inviável a realização de provas ou outras formas similares de avaliação. Esta se dará no dia a dia do aluno, levando em conta: assiduidade, pontualidade e material completo na pasta; participação construtiva em sala de aula e nas apresentações públicas - prontidão, envolvimento e seu real aproveitamento vocal e musical.'
$ws.Range("C20").Value = 'Sendo uma atividade prática e de grupo, fica inviável a realização de provas ou outras formas similares de avaliação. Esta se dará no dia a dia do aluno, levando em conta: assiduidade, pontualidade e material completo na pasta; participação construtiva em sala de aula e nas apresentações públicas - prontidão, envolvimento e seu real aproveitamento vocal e musical.'

$ws.Range("A21").Value = 'Bibliografia:'
$ws.Range("B21").Value = 'não tem'
$ws.Range("C21").Value = 'não tem'

$ws.Range("A22").Value = 'Requisitos:'
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()

$ws.Range("A23").ClearContents()
$ws.Range("B23").Value = "8800010 -  Canto Coral II  (Requisito)`n"
$ws.Range("C23").Value = "8800010 -  Canto Coral II  (Requisito)`n"

# --- Fix up row heights that changed as content moved between rows ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(17).EntireRow.AutoFit()  # back to default height
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).EntireRow.AutoFit()  # back to default height
$ws.Rows.Item(23).RowHeight = 30

# --- Row 24 (old trailing "Requisitos" detail row) is now gone entirely ---
$ws.Range("A24").EntireRow.Delete()
